$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '27.578.85'
Set-TextValue $ws.Range('E2') '  -1.16%  '

Set-TextValue $ws.Range('D3') '1.846.00'
Set-TextValue $ws.Range('E3') '  -0.90%  '

Set-TextValue $ws.Range('E4') '  -0.04%  '

Set-TextValue $ws.Range('D5') '314.32'
Set-TextValue $ws.Range('E5') '  -1.30%  '

Set-TextValue $ws.Range('D6') '1.003'
Set-TextValue $ws.Range('E6') '  -0.17%  '

Set-TextValue $ws.Range('D7') '0.4268'
Set-TextValue $ws.Range('E7') '  -2.21%  '

Set-TextValue $ws.Range('D8') '0.3655'
Set-TextValue $ws.Range('E8') '  -2.05%  '

Set-TextValue $ws.Range('D9') '44.75'
Set-TextValue $ws.Range('E9') '  +0.42%  '

Set-TextValue $ws.Range('D10') '0.07314'
Set-TextValue $ws.Range('E10') '  -2.35%  '

Set-TextValue $ws.Range('D11') '0.8851'
Set-TextValue $ws.Range('E11') '  -5.40%  '

Set-TextValue $ws.Range('D12') '20.85'
Set-TextValue $ws.Range('E12') '  -1.85%  '

Set-TextValue $ws.Range('D13') '1.909.49'
Set-TextValue $ws.Range('E13') '  -2.15%  '

Set-TextValue $ws.Range('D14') '5.359'
Set-TextValue $ws.Range('E14') '  -1.52%  '

Set-TextValue $ws.Range('D15') '6.552'
Set-TextValue $ws.Range('E15') '  -2.60%  '

Set-TextValue $ws.Range('D16') '0.06932'
Set-TextValue $ws.Range('E16') '  +0.47%  '

Set-TextValue $ws.Range('E17') '  -0.14%  '

Set-TextValue $ws.Range('D18') '78.82'
Set-TextValue $ws.Range('E18') '  -3.22%  '

Set-TextValue $ws.Range('D19') '0.000008889'
Set-TextValue $ws.Range('E19') '  -1.76%  '

Set-TextValue $ws.Range('E20') '  -0.01%  '

Set-TextValue $ws.Range('E21') '  -2.62%  '

Set-TextValue $ws.Range('D22') '27.595.39'
Set-TextValue $ws.Range('E22') '  -1.03%  '

Set-TextValue $ws.Range('D23') '4.990'
Set-TextValue $ws.Range('E23') '  -2.61%  '

Set-TextValue $ws.Range('D24') '10.67'
Set-TextValue $ws.Range('E24') '  -3.37%  '

Set-TextValue $ws.Range('D25') '2.122.18'
Set-TextValue $ws.Range('E25') '  -1.72%  '

Set-TextValue $ws.Range('D26') '1.962'
Set-TextValue $ws.Range('E26') '  -2.36%  '

Set-TextValue $ws.Range('D27') '153.91'
Set-TextValue $ws.Range('E27') '  -0.51%  '

Set-TextValue $ws.Range('D28') '18.98'
Set-TextValue $ws.Range('E28') '  +2.60%  '

Set-TextValue $ws.Range('D29') '121.91'
Set-TextValue $ws.Range('E29') '  +7.37%  '

Set-TextValue $ws.Range('D30') '5.242'
Set-TextValue $ws.Range('E30') '  -5.44%  '

Set-TextValue $ws.Range('D31') '1.919'
Set-TextValue $ws.Range('E31') '  +12.17%  '

Set-TextValue $ws.Range('D32') '0.08926'
Set-TextValue $ws.Range('E32') '  -1.09%  '

Set-TextValue $ws.Range('D33') '0.7626'
Set-TextValue $ws.Range('E33') '  -6.53%  '

Set-TextValue $ws.Range('D34') '4.575'
Set-TextValue $ws.Range('E34') '  -5.21%  '

Set-TextValue $ws.Range('D35') '2.967'
Set-TextValue $ws.Range('E35') '  -0.21%  '

Set-TextValue $ws.Range('D36') '1.100'
Set-TextValue $ws.Range('E36') '  -6.52%  '

Set-TextValue $ws.Range('E37') '  -0.24%  '

Set-TextValue $ws.Range('B38') 'TrustWalletToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D38') '1.097'
Set-TextValue $ws.Range('E38') '  -2.12%  '

Set-TextValue $ws.Range('B39') 'Hedera'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D39') '0.05380'
Set-TextValue $ws.Range('E39') '  -2.62%  '

Set-TextValue $ws.Range('D40') '0.01948'
Set-TextValue $ws.Range('E40') '  -1.46%  '

Set-TextValue $ws.Range('D41') '2.807'
Set-TextValue $ws.Range('E41') '  -4.82%  '

Set-TextValue $ws.Range('D42') '6.920'
Set-TextValue $ws.Range('E42') '  -1.46%  '

Set-TextValue $ws.Range('D43') '0.5115'
Set-TextValue $ws.Range('E43') '  -2.93%  '

Set-TextValue $ws.Range('D44') '0.1656'
Set-TextValue $ws.Range('E44') '  -2.59%  '

Set-TextValue $ws.Range('D45') '8.273'
Set-TextValue $ws.Range('E45') '  -5.98%  '

Set-TextValue $ws.Range('D46') '0.06580'
Set-TextValue $ws.Range('E46') '  -2.53%  '

Set-TextValue $ws.Range('B47') 'EnergySwap'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D47') '10.45'
Set-TextValue $ws.Range('E47') '  -1.34%  '

Set-TextValue $ws.Range('B48') 'Decentraland'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range('D48') '0.4762'
Set-TextValue $ws.Range('E48') '  -2.71%  '

Set-TextValue $ws.Range('D49') '104.31'
Set-TextValue $ws.Range('E49') '  -3.28%  '

Set-TextValue $ws.Range('D50') '1.002'
Set-TextValue $ws.Range('E50') '  -0.20%  '

Set-TextValue $ws.Range('E51') '  -2.71%  '
